$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date formatting (style s="2") from A269 down through the new rows
$ws.Range("A269").Copy() | Out-Null
$ws.Range("A270:A301").PasteSpecial(-4122) | Out-Null

$rows = @(
    @(270,44344,1,9,35.22642764883166),
    @(271,44345,1,10,39.14047516536851),
    @(272,44346,1,10,39.14047516536851),
    @(273,44347,4,14,54.79666523151591),
    @(274,44348,0,14,54.79666523151591),
    @(275,44349,0,13,50.88261771497906),
    @(276,44350,1,8,31.31238013229481),
    @(277,44351,1,8,31.31238013229481),
    @(278,44352,0,7,27.39833261575795),
    @(279,44353,1,7,27.39833261575795),
    @(280,44354,2,5,19.57023758268425),
    @(281,44355,1,6,23.4842850992211),
    @(282,44356,0,6,23.4842850992211),
    @(283,44357,1,6,23.4842850992211),
    @(284,44358,2,7,27.39833261575795),
    @(285,44359,0,7,27.39833261575795),
    @(286,44360,0,6,23.4842850992211),
    @(287,44361,2,6,23.4842850992211),
    @(288,44362,0,5,19.57023758268425),
    @(289,44363,1,6,23.4842850992211),
    @(290,44364,3,8,31.31238013229481),
    @(291,44365,0,6,23.4842850992211),
    @(292,44366,3,9,35.22642764883166),
    @(293,44367,0,9,35.22642764883166),
    @(294,44368,1,8,31.31238013229481),
    @(295,44369,1,9,35.22642764883166),
    @(296,44370,0,8,31.31238013229481),
    @(297,44371,1,6,23.4842850992211),
    @(298,44372,0,6,23.4842850992211),
    @(299,44373,0,3,11.74214254961055),
    @(300,44374,2,5,19.57023758268425),
    @(301,44375,1,5,19.57023758268425)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row,1).Value = $r[1]
    $ws.Cells.Item($row,2).Value = $r[2]
    $ws.Cells.Item($row,3).Value = $r[3]
    $ws.Cells.Item($row,4).Value = $r[4]
}

Write-Host "done"
